$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the BabyDogeCoin insertion
# shifting Cronos/EnergySwap/Algorand down one row, dropping Mantle).
# Cells whose new text looks like a plain number (e.g. "4.20") are forced
# to Text format first so Excel keeps the exact string instead of
# collapsing it to a numeric value (which would drop trailing zeros, etc).

$ws.Range('D2').Value = '26.657.91'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.643.06'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.08'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0627'
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.29'
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.872.15'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '1.665.48'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.20'
$ws.Range('E14').Value = '  +2.65%  '
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.51'
$ws.Range('E16').Value = '  +3.22%  '
$ws.Range('D17').Value = '26.723.71'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '0.0₃0748'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.75'
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.31'
$ws.Range('E22').Value = '  +2.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.54'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.19'
$ws.Range('E24').Value = '  +12.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.84'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('E28').Value = '  +4.55%  '
$ws.Range('E29').Value = '  +1.50%  '
$ws.Range('E30').Value = '  +2.40%  '
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('D34').Value = '1.273.93'
$ws.Range('E34').Value = '  +4.64%  '
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0181'
$ws.Range('E36').Value = '  +5.99%  '
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.534'
$ws.Range('E38').Value = '  +6.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.829'
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E41').Value = '  +2.66%  '
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').Value = '1.782.45'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.49'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.93'
$ws.Range('E46').Value = '  +8.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.61'
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0516'
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.82'
$ws.Range('E50').Value = '  +2.32%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0972'
$ws.Range('E51').Value = '  +3.18%  '
